# "Generate Report for Archive" — refresh the localization-status report:
#   * the two in-flight hand-off rows moved from "Ready for handoff" to
#     "In Translation" on every sheet that tracks them
#   * the now-shorter status text means the status column(s) no longer
#     need to be as wide, so re-fit them to the new content

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet: status is mirrored per-locale in columns E (zh-cn) and F (de-de)
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# --- Per-locale detail sheets: status lives in column C
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# --- Re-fit the status columns now that the text is shorter than
# "Ready for handoff". (12.5 is the ColumnWidth input that lands on the
# nearest achievable pixel-grid width to the new autofit target.)
$newColumnWidth = 12.5
$overview.Columns.Item(5).ColumnWidth = $newColumnWidth
$overview.Columns.Item(6).ColumnWidth = $newColumnWidth
$zhcn.Columns.Item(3).ColumnWidth = $newColumnWidth
$dede.Columns.Item(3).ColumnWidth = $newColumnWidth
